$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'36.140.21"
$ws.Range("E2").Value = "'  +2.23%  "

# Row 3
$ws.Range("D3").Value = "'2.012.42"
$ws.Range("E3").Value = "'  +6.16%  "

# Row 4
$ws.Range("E4").Value = "'  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'245.31"
$ws.Range("E5").Value = "'  -0.57%  "

# Row 6
$ws.Range("E6").Value = "'  -4.96%  "

# Row 7
$ws.Range("E7").Value = "'  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'44.79"
$ws.Range("E8").Value = "'  +4.26%  "

# Row 9
$ws.Range("D9").Value = "'58.36"
$ws.Range("E9").Value = "'  +4.41%  "

# Row 10
$ws.Range("D10").Value = "'0.360"
$ws.Range("E10").Value = "'  +0.63%  "

# Row 11
$ws.Range("D11").Value = "'0.0712"
$ws.Range("E11").Value = "'  -4.91%  "

# Row 12
$ws.Range("D12").Value = "'0.0986"
$ws.Range("E12").Value = "'  +0.41%  "

# Row 13
$ws.Range("D13").Value = "'14.55"
$ws.Range("E13").Value = "'  +4.24%  "

# Row 14
$ws.Range("D14").Value = "'2.305.44"
$ws.Range("E14").Value = "'  +6.15%  "

# Row 15
$ws.Range("D15").Value = "'0.800"
$ws.Range("E15").Value = "'  +0.74%  "

# Row 16
$ws.Range("D16").Value = "'2.022.11"
$ws.Range("E16").Value = "'  +6.72%  "

# Row 17
$ws.Range("D17").Value = "'4.87"
$ws.Range("E17").Value = "'  -2.77%  "

# Row 18
$ws.Range("D18").Value = "'36.259.95"
$ws.Range("E18").Value = "'  +2.37%  "

# Row 19
$ws.Range("D19").Value = "'70.69"
$ws.Range("E19").Value = "'  -3.83%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0815"
$ws.Range("E20").Value = "'  -1.56%  "

# Row 21
$ws.Range("D21").Value = "'12.90"
$ws.Range("E21").Value = "'  -0.31%  "

# Row 22
$ws.Range("D22").Value = "'233.79"
$ws.Range("E22").Value = "'  -4.42%  "

# Row 23
$ws.Range("D23").Value = "'4.90"
$ws.Range("E23").Value = "'  -6.39%  "

# Row 24
$ws.Range("E24").Value = "'  +0.05%  "

# Row 25
$ws.Range("E25").Value = "'  -8.13%  "

# Row 26
$ws.Range("D26").Value = "'161.95"
$ws.Range("E26").Value = "'  -2.75%  "

# Row 27
$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.65"
$ws.Range("E27").Value = "'  +7.22%  "

# Row 28
$ws.Range("B28").Value = "'PancakeSwap"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'1.97"
$ws.Range("E28").Value = "'  -9.31%  "

# Row 29
$ws.Range("D29").Value = "'8.48"
$ws.Range("E29").Value = "'  -0.80%  "

# Row 30
$ws.Range("E30").Value = "'  -4.87%  "

# Row 31
$ws.Range("D31").Value = "'4.36"
$ws.Range("E31").Value = "'  +0.29%  "

# Row 32
$ws.Range("D32").Value = "'20.97"
$ws.Range("E32").Value = "'  +53.80%  "

# Row 33
$ws.Range("E33").Value = "'  -3.09%  "

# Row 34
$ws.Range("E34").Value = "'  +0.06%  "

# Row 35
$ws.Range("E35").Value = "'  -0.68%  "

# Row 36
$ws.Range("D36").Value = "'3.99"
$ws.Range("E36").Value = "'  -5.67%  "

# Row 37
$ws.Range("D37").Value = "'0.0809"
$ws.Range("E37").Value = "'  +12.77%  "

# Row 38
$ws.Range("D38").Value = "'2.13"
$ws.Range("E38").Value = "'  +8.35%  "

# Row 39
$ws.Range("D39").Value = "'0.833"
$ws.Range("E39").Value = "'  -2.29%  "

# Row 40
$ws.Range("D40").Value = "'1.34"
$ws.Range("E40").Value = "'  -8.70%  "

# Row 41
$ws.Range("B41").Value = "'VeChain"
$ws.Range("C41").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0214"
$ws.Range("E41").Value = "'  -4.96%  "

# Row 42
$ws.Range("B42").Value = "'Aave"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'95.75"
$ws.Range("E42").Value = "'  -3.36%  "

# Row 43
$ws.Range("E43").Value = "'  -0.25%  "

# Row 44
$ws.Range("D44").Value = "'2.74"
$ws.Range("E44").Value = "'  +13.70%  "

# Row 45
$ws.Range("E45").Value = "'  -6.32%  "

# Row 46
$ws.Range("D46").Value = "'1.314.16"
$ws.Range("E46").Value = "'  -1.81%  "

# Row 47
$ws.Range("D47").Value = "'0.0808"
$ws.Range("E47").Value = "'  -0.08%  "

# Row 48
$ws.Range("B48").Value = "'RocketPoolETH"
$ws.Range("C48").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'2.231.28"
$ws.Range("E48").Value = "'  +7.64%  "

# Row 49
$ws.Range("B49").Value = "'MXToken"
$ws.Range("C49").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.77"
$ws.Range("E49").Value = "'  +1.15%  "

# Row 50
$ws.Range("D50").Value = "'2.20"
$ws.Range("E50").Value = "'  -7.53%  "

# Row 51
$ws.Range("D51").Value = "'3.76"
$ws.Range("E51").Value = "'  +12.48%  "
